$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02542733333333333
$ws.Range("H2").Value = 0.076282
$ws.Range("I2").Value = 0.3241752404264994
$ws.Range("J2").Value = 0.3241752404264994
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 0.1603408481906667
$ws.Range("R2").Value = 1.443067633716
$ws.Range("S2").Value = 0.004396304182554078
$ws.Range("T2").Value = 0.004396304182554078

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02542733333333333
$ws.Range("H3").Value = 0.076282
$ws.Range("I3").Value = 0.3241752404264994
$ws.Range("J3").Value = 0.3241752404264994
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 4.641293930583778
$ws.Range("R3").Value = 41.771645375254
$ws.Range("S3").Value = 0.1272572781654778
$ws.Range("T3").Value = 0.1272572781654778

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02542733333333333
$ws.Range("H4").Value = 0.076282
$ws.Range("I4").Value = 0.3241752404264994
$ws.Range("J4").Value = 0.3241752404264994
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 3.239350550275333
$ws.Range("R4").Value = 29.154154952478
$ws.Range("S4").Value = 0.0888181055148196
$ws.Range("T4").Value = 0.08881810551481961

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.02542733333333333
$ws.Range("H5").Value = 0.076282
$ws.Range("I5").Value = 0.3241752404264994
$ws.Range("J5").Value = 0.3241752404264994
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 0.4939005565697778
$ws.Range("R5").Value = 4.445105009128
$ws.Range("S5").Value = 0.01354200821010683
$ws.Range("T5").Value = 0.01354200821010683

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02542733333333333
$ws.Range("H6").Value = 0.076282
$ws.Range("I6").Value = 0.3241752404264994
$ws.Range("J6").Value = 0.3241752404264994
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 3.288348097748889
$ws.Range("R6").Value = 29.59513287974
$ws.Range("S6").Value = 0.09016154435354108
$ws.Range("T6").Value = 0.09016154435354108

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.05300966666666667
$ws.Range("H7").Value = 0.159029
$ws.Range("I7").Value = 0.6758247595735006
$ws.Range("J7").Value = 0.6758247595735006
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 0.3342707945113333
$ws.Range("R7").Value = 3.008437150602
$ws.Range("S7").Value = 0.00916520093662191
$ws.Range("T7").Value = 0.009165200936621909

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.05300966666666667
$ws.Range("H8").Value = 0.159029
$ws.Range("I8").Value = 0.6758247595735006
$ws.Range("J8").Value = 0.6758247595735006
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("Q8").Value = 9.67594363659589
$ws.Range("R8").Value = 87.083492729363
$ws.Range("S8").Value = 0.2652997783143832
$ws.Range("T8").Value = 0.2652997783143831

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.05300966666666667
$ws.Range("H9").Value = 0.159029
$ws.Range("I9").Value = 0.6758247595735006
$ws.Range("J9").Value = 0.6758247595735006
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 6.753240327465667
$ws.Range("R9").Value = 60.779162947191
$ws.Range("S9").Value = 0.1851636624880869
$ws.Range("T9").Value = 0.1851636624880869

# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.05300966666666667
$ws.Range("H10").Value = 0.159029
$ws.Range("I10").Value = 0.6758247595735006
$ws.Range("J10").Value = 0.6758247595735006
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 1.029659836012889
$ws.Range("R10").Value = 9.266938524116
$ws.Range("S10").Value = 0.02823171945734354
$ws.Range("T10").Value = 0.02823171945734353

# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.05300966666666667
$ws.Range("H11").Value = 0.159029
$ws.Range("I11").Value = 0.6758247595735006
$ws.Range("J11").Value = 0.6758247595735006
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 6.855388029114446
$ws.Range("R11").Value = 61.69849226203001
$ws.Range("S11").Value = 0.1879643983770652
$ws.Range("T11").Value = 0.1879643983770651
